$d = $word.ActiveDocument

# 1. "Build attractive, intuitive graphical user interfaces"
#    -> split into two runs with an empty _GoBack bookmark between them.
$r = $d.Content
$r.Find.Execute("Build attractive, intuitive graphical ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $splitPoint)

# 2. "Justify the choice to use an object-oriented approach to developing software"
#    -> "Justify the choice an object-oriented approach to developing software"
$d.Content.Find.Execute("Justify the choice to use an object-oriented approach to developing software", $true, $false, $false, $false, $false, $true, 1, $false, "Justify the choice an object-oriented approach to developing software", 2)

# 3. Merge "Homework will be assigned weekly to he" + bookmark + "lp you keep..."
#    into a single run, removing the old _GoBack bookmark in the process.
$d.Content.Find.Execute("Homework will be assigned weekly to help you keep your skills sharp.", $true, $false, $false, $false, $false, $true, 1, $false, "Homework will be assigned weekly to help you keep your skills sharp.", 2)
